$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3669.9
$ws.Range("I40").Value = 3628.1428
$ws.Range("J40").Value = 3767.3333
$ws.Range("K40").Value = 3628.1428
$ws.Range("L40").Value = 3767.3333
$ws.Range("M40").Value = -3453.1428
$ws.Range("N40").Value = -4117.3333
$ws.Range("H70").Value = 5017.391
$ws.Range("I70").Value = 1764.5385
$ws.Range("K70").Value = 5293.6155
$ws.Range("M70").Value = -5023.6155
$ws.Range("H73").Value = 5017.391
$ws.Range("I73").Value = 1764.5385
$ws.Range("K73").Value = 5293.6155
$ws.Range("M73").Value = -4357.6155
$ws.Range("H88").Value = 4184.4287
$ws.Range("J88").Value = 5497
$ws.Range("L88").Value = 5497
$ws.Range("N88").Value = -6309
$ws.Range("H91").Value = 4184.4287
$ws.Range("J91").Value = 5497
$ws.Range("L91").Value = 5497
$ws.Range("N91").Value = -8305
$ws.Range("H106").Value = 87180.234
$ws.Range("I106").Value = 110332.5
$ws.Range("K106").Value = 110332.5
$ws.Range("M106").Value = -109701.5
$ws.Range("H132").Value = 1623.0294
$ws.Range("J132").Value = 1283.2
$ws.Range("L132").Value = 3849.6
$ws.Range("N132").Value = -8909.6
$ws.Range("H134").Value = 108831.664
$ws.Range("J134").Value = 108831.664
$ws.Range("L134").Value = 108831.664
$ws.Range("N134").Value = -118971.664
$ws.Range("H137").Value = 130853.29
$ws.Range("I137").Value = 2117.2144
$ws.Range("J137").Value = 671544.8
$ws.Range("K137").Value = 6351.6432
$ws.Range("L137").Value = 2014634.4
$ws.Range("M137").Value = -3801.6432
$ws.Range("N137").Value = -2019734.4
$ws.Range("H138").Value = 2475
$ws.Range("I138").Value = 1824.1052
$ws.Range("J138").Value = 4241.7144
$ws.Range("K138").Value = 5472.3156
$ws.Range("L138").Value = 12725.1432
$ws.Range("M138").Value = -332.3155999999999
$ws.Range("N138").Value = -23005.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1309.5
$ws.Range("I45").Value = 1157.9375
$ws.Range("K45").Value = 1157.9375
$ws.Range("M45").Value = -780.9375
$ws.Range("H61").Value = 14374.213
$ws.Range("I61").Value = 2327.5925
$ws.Range("J61").Value = 39394.117
$ws.Range("K61").Value = 2327.5925
$ws.Range("L61").Value = 39394.117
$ws.Range("M61").Value = -2115.5925
$ws.Range("N61").Value = -39818.117
$ws.Range("H95").Value = 23534.5
$ws.Range("J95").Value = 23534.5
$ws.Range("L95").Value = 23534.5
$ws.Range("N95").Value = -29026.5
$ws.Range("H110").Value = 96067.52
$ws.Range("I110").Value = 96067.52
$ws.Range("K110").Value = 96067.52
$ws.Range("M110").Value = -94022.52
$ws.Range("H132").Value = 7412.75
$ws.Range("I132").Value = 8117.136
$ws.Range("K132").Value = 24351.408
$ws.Range("M132").Value = -21821.408
$ws.Range("H136").Value = 14374.213
$ws.Range("I136").Value = 2327.5925
$ws.Range("J136").Value = 39394.117
$ws.Range("K136").Value = 6982.7775
$ws.Range("L136").Value = 118182.351
$ws.Range("M136").Value = -4432.7775
$ws.Range("N136").Value = -123282.351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 15499.5
$ws.Range("I18").Value = 1000
$ws.Range("J18").Value = 29999
$ws.Range("K18").Value = 1000
$ws.Range("L18").Value = 29999
$ws.Range("M18").Value = -471
$ws.Range("N18").Value = -31057
$ws.Range("H99").Value = 1681.3077
$ws.Range("I99").Value = 1860.7
$ws.Range("J99").Value = 1083.3334
$ws.Range("K99").Value = 1860.7
$ws.Range("L99").Value = 1083.3334
$ws.Range("M99").Value = -362.7
$ws.Range("N99").Value = -4079.3334
$ws.Range("H117").Value = 71666.664
$ws.Range("J117").Value = 71666.664
$ws.Range("L117").Value = 71666.664
$ws.Range("N117").Value = -80844.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2236.9136
$ws.Range("I31").Value = 1704.5128
$ws.Range("J31").Value = 2731.2856
$ws.Range("K31").Value = 1704.5128
$ws.Range("L31").Value = 2731.2856
$ws.Range("M31").Value = -1409.5128
$ws.Range("N31").Value = -3321.2856
$ws.Range("H34").Value = 2236.9136
$ws.Range("I34").Value = 1704.5128
$ws.Range("J34").Value = 2731.2856
$ws.Range("K34").Value = 1704.5128
$ws.Range("L34").Value = 2731.2856
$ws.Range("M34").Value = -1502.5128
$ws.Range("N34").Value = -3135.2856
$ws.Range("H132").Value = 8108.1763
$ws.Range("I132").Value = 2988.5715
$ws.Range("K132").Value = 8965.7145
$ws.Range("M132").Value = -6435.7145
$ws.Range("H134").Value = 2472.6072
$ws.Range("I134").Value = 1974.7273
$ws.Range("K134").Value = 5924.1819
$ws.Range("M134").Value = -3389.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 4635.273
$ws.Range("I69").Value = 2995
$ws.Range("K69").Value = 8985
$ws.Range("M69").Value = -8174
$ws.Range("H72").Value = 4635.273
$ws.Range("I72").Value = 2995
$ws.Range("K72").Value = 26955
$ws.Range("M72").Value = -22899
$ws.Range("H119").Value = 10400
$ws.Range("I119").Value = 10400
$ws.Range("K119").Value = 31200
$ws.Range("M119").Value = -26362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 39139.543
$ws.Range("I122").Value = 48810.05
$ws.Range("J122").Value = 2391.6
$ws.Range("K122").Value = 146430.15
$ws.Range("L122").Value = 7174.799999999999
$ws.Range("M122").Value = -143980.15
$ws.Range("N122").Value = -12074.8
$ws.Range("H132").Value = 2434.25
$ws.Range("I132").Value = 2701.4524
$ws.Range("K132").Value = 8104.3572
$ws.Range("M132").Value = -5574.3572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3492.5366
$ws.Range("I22").Value = 6540.067
$ws.Range("J22").Value = 1734.3462
$ws.Range("K22").Value = 6540.067
$ws.Range("L22").Value = 1734.3462
$ws.Range("M22").Value = -6245.067
$ws.Range("N22").Value = -2324.3462
$ws.Range("H27").Value = 3492.5366
$ws.Range("I27").Value = 6540.067
$ws.Range("J27").Value = 1734.3462
$ws.Range("K27").Value = 6540.067
$ws.Range("L27").Value = 1734.3462
$ws.Range("M27").Value = -6433.067
$ws.Range("N27").Value = -1948.3462
$ws.Range("H40").Value = 5099.727
$ws.Range("I40").Value = 4533.409
$ws.Range("J40").Value = 6232.364
$ws.Range("K40").Value = 4533.409
$ws.Range("L40").Value = 6232.364
$ws.Range("M40").Value = -4397.409
$ws.Range("N40").Value = -6504.364
$ws.Range("H50").Value = 32054.666
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H55").Value = 1133.0344
$ws.Range("I55").Value = 397.17392
$ws.Range("K55").Value = 397.17392
$ws.Range("M55").Value = -224.17392
$ws.Range("H56").Value = 17050
$ws.Range("I56").Value = 12200
$ws.Range("J56").Value = 21900
$ws.Range("K56").Value = 12200
$ws.Range("L56").Value = 21900
$ws.Range("M56").Value = -11509
$ws.Range("N56").Value = -23282
$ws.Range("H57").Value = 30050
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H61").Value = 27028380
$ws.Range("H113").Value = 27028380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 890.8929000000001
$ws.Range("I113").Value = 921.0769
$ws.Range("J113").Value = 2763.8076
$ws.Range("K113").Value = 2763.2307
$ws.Range("M113").Value = -593.2307000000001

